$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("B2").Value = 98
$ws.Range("D2").Value = 78.358873500000001

# Remove row 3 entirely
$ws.Rows.Item(3).Delete()
